$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '37.975.99'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +2.51%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.054.45'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +2.13%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -1.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '229.68'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +1.63%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.615'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +1.35%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '57.86'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +6.18%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.386'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +1.89%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0803'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +2.12%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -1.14%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.363.68'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +2.20%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.57'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +2.71%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.58'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +1.86%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.751'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +1.61%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.26'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +3.01%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.060.60'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +2.19%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '37.926.35'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +2.58%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.18'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +1.58%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '69.70'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +1.35%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0826'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +1.11%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '224.09'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.31%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.06%  '
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.68%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.25'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +2.49%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.26'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.96%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '165.47'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.41%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.134'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +7.07%  '
$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = 'EthereumClassic'
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.01'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +1.79%  '
$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = 'ImmutableX'
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.37'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +1.10%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.118'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +1.72%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.53'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.74%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.59'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +3.86%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0612'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.19%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.37'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +1.56%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.01'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +11.79%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.30'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +5.34%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -0.20%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '98.41'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +3.80%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0217'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.74%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.476.56'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.29%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0941'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +2.31%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.86'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +3.97%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '16.73'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +1.43%  '
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.19%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.13'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +18.36%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +1.14%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +1.89%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.08'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -1.97%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.253.53'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +2.39%  '
